# Rename "Sheet1" to "SBS".
# Excel automatically updates all formulas that reference the sheet
# (e.g. COUNTIF(Sheet1!$B$2:$B$73, ...) -> COUNTIF(SBS!$B$2:$B$73, ...))
# when the worksheet's .Name property is changed.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Name = "SBS"
